$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# Rows 1-3: memory figures became "0M" placeholders
$tbl.Cell(1,1).Range.Text = "0M"
$tbl.Cell(2,1).Range.Text = "0M"
$tbl.Cell(3,1).Range.Text = "0M"

# Row 4: GC count updated
$tbl.Cell(4,1).Range.Text = "32"

# Rows 6,7,8,10: refreshed pause-time figures (row 5, 9, 11 unchanged)
$tbl.Cell(6,1).Range.Text = "0.00005"
$tbl.Cell(7,1).Range.Text = "0.00004"
$tbl.Cell(8,1).Range.Text = "0.00001"
$tbl.Cell(10,1).Range.Text = "0.00004"

# Row 12: total pause time updated
$tbl.Cell(12,1).Range.Text = "0.00117"

# Last three rows: collapse the multi-column tab-separated stats down to a
# single summary value each.
$tbl.Cell(44,1).Range.Text = "100"
$tbl.Cell(45,1).Range.Text = "0"
$tbl.Cell(46,1).Range.Text = "97"
